# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# The account-statement table (rows 16-56) lists one row per overdue period
# ("Periodo Mora", column E), with its "Valor Mora" (column F) and
# "Salario Basico" (column G). The sheet is refreshed with new data:
#   - the periods are re-sorted from newest-first to oldest-first
#     (chronological order, 1611 .. 2003 instead of 2003 .. 1611)
#   - "Valor Mora" (F) is updated per period
#   - "Salario Basico" (G) is updated to the new uniform value

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$periods = @( `
    "1611","1612", `
    "1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712", `
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812", `
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912", `
    "2001","2002","2003" `
)

$firstRow = 16
$newSalarioBasico = 781242

for ($i = 0; $i -lt $periods.Count; $i++) {
    $row = $firstRow + $i
    $period = $periods[$i]

    if ($row -le 37) {
        $valorMora = 27578
    } else {
        $valorMora = 31249
    }

    $ws.Range("E$row").Value = $period
    $ws.Range("F$row").Value = $valorMora
    $ws.Range("G$row").Value = $newSalarioBasico
}
